$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new value is a plain decimal number must be
# forced to Text format first, matching the source data (inline strings,
# e.g. "1.003", "0.1050") -- otherwise Excel auto-converts them to numbers
# and mangles formatting (trailing zeros, leading zeros, etc).
$textCells = @("D4","D5","D7","D8","D9","D10","D13","D14","D16","D17","D18","D19","D20","D21","D24","D25","D26","D27","D28","D29","D30","D32","D33","D34","D35","D36","D37","D38","D39","D41","D43","D44","D45","D46","D47","D48","D49","D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "30.396.82"
$ws.Range("E2").Value = "  +2.25%  "
$ws.Range("D3").Value = "2.095.01"
$ws.Range("E3").Value = "  -0.03%  "
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  -0.65%  "
$ws.Range("D5").Value = "342.88"
$ws.Range("E5").Value = "  -0.13%  "
$ws.Range("E6").Value = "  -0.66%  "
$ws.Range("D7").Value = "0.5247"
$ws.Range("E7").Value = "  +1.60%  "
$ws.Range("D8").Value = "0.4428"
$ws.Range("E8").Value = "  +1.14%  "
$ws.Range("D9").Value = "54.49"
$ws.Range("E9").Value = "  +3.36%  "
$ws.Range("D10").Value = "0.09369"
$ws.Range("E10").Value = "  +0.99%  "
$ws.Range("E11").Value = "  +0.58%  "
$ws.Range("E12").Value = "  -0.36%  "
$ws.Range("D13").Value = "8.588"
$ws.Range("E13").Value = "  +3.88%  "
$ws.Range("D14").Value = "6.923"
$ws.Range("E14").Value = "  +2.60%  "
$ws.Range("D15").Value = "2.064.07"
$ws.Range("E15").Value = "  -1.97%  "
$ws.Range("D16").Value = "101.49"
$ws.Range("E16").Value = "  +1.97%  "
$ws.Range("D17").Value = "0.00001159"
$ws.Range("E17").Value = "  +0.48%  "
$ws.Range("D18").Value = "1.003"
$ws.Range("E18").Value = "  -0.69%  "
$ws.Range("D19").Value = "21.19"
$ws.Range("E19").Value = "  +1.92%  "
$ws.Range("D20").Value = "0.06676"
$ws.Range("E20").Value = "  +0.49%  "
$ws.Range("D21").Value = "6.330"
$ws.Range("E21").Value = "  +2.26%  "
$ws.Range("E22").Value = "  -0.58%  "
$ws.Range("D23").Value = "30.434.53"
$ws.Range("E23").Value = "  +2.28%  "
$ws.Range("D24").Value = "12.56"
$ws.Range("E24").Value = "  +0.58%  "
$ws.Range("D25").Value = "2.311"
$ws.Range("E25").Value = "  -0.25%  "
$ws.Range("D26").Value = "21.86"
$ws.Range("E26").Value = "  -0.30%  "
$ws.Range("D27").Value = "162.94"
$ws.Range("E27").Value = "  +0.99%  "
$ws.Range("D28").Value = "6.799"
$ws.Range("E28").Value = "  +8.56%  "
$ws.Range("D29").Value = "2.512"
$ws.Range("E29").Value = "  +0.04%  "
$ws.Range("D30").Value = "133.52"
$ws.Range("E30").Value = "  +0.39%  "
$ws.Range("E31").Value = "  +0.33%  "
$ws.Range("B32").Value = "Stellar"
$ws.Range("C32").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D32").Value = "0.1050"
$ws.Range("E32").Value = "  +0.15%  "
$ws.Range("B33").Value = "ARBITRUM"
$ws.Range("C33").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D33").Value = "1.660"
$ws.Range("E33").Value = "  +0.56%  "
$ws.Range("D34").Value = "6.272"
$ws.Range("E34").Value = "  +1.83%  "
$ws.Range("D35").Value = "3.878"
$ws.Range("E35").Value = "  -1.52%  "
$ws.Range("D36").Value = "10.17"
$ws.Range("E36").Value = "  -0.32%  "
$ws.Range("D37").Value = "0.02632"
$ws.Range("E37").Value = "  +2.18%  "
$ws.Range("D38").Value = "0.06822"
$ws.Range("E38").Value = "  +1.61%  "
$ws.Range("D39").Value = "0.7007"
$ws.Range("E39").Value = "  +1.59%  "
$ws.Range("E40").Value = "  +1.14%  "
$ws.Range("D41").Value = "1.342"
$ws.Range("E41").Value = "  +1.88%  "
$ws.Range("E42").Value = "  -0.18%  "
$ws.Range("D43").Value = "0.6843"
$ws.Range("E43").Value = "  +1.33%  "
$ws.Range("D44").Value = "14.40"
$ws.Range("E44").Value = "  +0.33%  "
$ws.Range("D45").Value = "2.346"
$ws.Range("E45").Value = "  +1.15%  "
$ws.Range("D46").Value = "1.001"
$ws.Range("E46").Value = "  -0.59%  "
$ws.Range("D47").Value = "1.384"
$ws.Range("E47").Value = "  +19.25%  "
$ws.Range("D48").Value = "3.634"
$ws.Range("E48").Value = "  +0.68%  "
$ws.Range("D49").Value = "0.00000000353"
$ws.Range("E49").Value = "  -1.68%  "
$ws.Range("D50").Value = "1.231"
$ws.Range("E50").Value = "  +9.87%  "
$ws.Range("E51").Value = "  -0.06%  "
